$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor name (Prakruti Sinha) - value field next to "Supervisor Name:" label
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor sign-off block: initials and sign-off date (28/02/2014),
# matching the style already used for the employee sign-off row above (D25).
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41698
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H27").Select()
